# Updating documents from invoices for December
# - "Hours by consultant" and "Cost" sheets: fill in the December 2024 (row 11)
#   actuals for TSB / MB / CGH / NJJ that were billed this month.
# - "Sheet1" (Northtech invoice breakdown): retitle for December, bump the
#   CGH/TSB hour counts to match, and add the new "Other" column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Hours by consultant": December 2024 is row 11 (A11 = "2024/12").
# Columns: B=TSB, C=MB, D=CGH, E=NJJ, H=TOTAL (formula), I=INVOICE TOTAL
# ---------------------------------------------------------------------
$hoursWs = $wb.Worksheets.Item("Hours by consultant")
$hoursWs.Range("B11").Value = 74
$hoursWs.Range("C11").Value = 0
$hoursWs.Range("D11").Value = 11
$hoursWs.Range("E11").Value = 0
$hoursWs.Range("I11").Value = 85

# ---------------------------------------------------------------------
# "Cost": same December 2024 row (row 11), in DKK.
# ---------------------------------------------------------------------
$costWs = $wb.Worksheets.Item("Cost")
$costWs.Range("B11").Value = 85040.8
$costWs.Range("C11").Value = 0
$costWs.Range("D11").Value = 15171.75
$costWs.Range("E11").Value = 0
$costWs.Range("I11").Value = 100212.55

# ---------------------------------------------------------------------
# "Sheet1": Northtech invoice detail for the month.
# ---------------------------------------------------------------------
$northWs = $wb.Worksheets.Item("Sheet1")

# Retitle the banner from October to December.
$northWs.Range("A1").Value = "Northtech Consultant Hours for: December 2024"

# Updated hour counts for this invoice (CGH row 4, TSB row 6).
$northWs.Range("B4").Value = 11
$northWs.Range("B6").Value = 74

# New "Other" work-package column, matching the formatting of the
# neighbouring "WP5B" column (J).
$northWs.Range("J3").Copy()
$northWs.Range("K3").PasteSpecial(-4122)
$northWs.Range("K3").Value = "Other"

$northWs.Range("J4:J7").Copy()
$northWs.Range("K4:K7").PasteSpecial(-4122)

$wb.Application.CutCopyMode = $false
